$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 347
$ws.Range("F2").Value = 45597
$ws.Range("G2").Value = 30468
$ws.Range("H2").Value = 45658
$ws.Range("E3").Value = 30407
$ws.Range("F3").Value = 45597
$ws.Range("E4").Value = 30407
$ws.Range("F4").Value = 45597
$ws.Range("G4").Value = 30468
$ws.Range("H4").Value = 45658
$ws.Range("E5").Value = 30407
$ws.Range("F5").Value = 45597
$ws.Range("G5").Value = 30468
$ws.Range("H5").Value = 45658
$ws.Range("C6").Value = 455
$ws.Range("F6").Value = 45566
$ws.Range("G6").Value = 30468
$ws.Range("H6").Value = 45658
$ws.Range("E7").Value = 30376
$ws.Range("F7").Value = 45566
$ws.Range("G7").Value = 30468
$ws.Range("H7").Value = 45658
$ws.Range("D8").Value = 423
$ws.Range("E8").Value = 30376
$ws.Range("F8").Value = 45566
$ws.Range("H8").Value = 45658
$ws.Range("E9").Value = 30407
$ws.Range("F9").Value = 45597
$ws.Range("G9").Value = 30468
$ws.Range("H9").Value = 45658
$ws.Range("D10").Value = 496
$ws.Range("E10").Value = 30407
$ws.Range("F10").Value = 45597
$ws.Range("H10").Value = 45658
$ws.Range("E11").Value = 30376
$ws.Range("F11").Value = 45566
$ws.Range("G11").Value = 30468
$ws.Range("H11").Value = 45658
$ws.Range("C12").Value = 384
$ws.Range("D12").Value = 365
$ws.Range("F12").Value = 45597
$ws.Range("H12").Value = 45658
$ws.Range("C13").Value = 480
$ws.Range("F13").Value = 45597
$ws.Range("G13").Value = 30468
$ws.Range("H13").Value = 45658
$ws.Range("C14").Value = 435
$ws.Range("D14").Value = 409
$ws.Range("F14").Value = 45566
$ws.Range("H14").Value = 45658
$ws.Range("G15").Value = 30438
$ws.Range("H15").Value = 45658
$ws.Range("C16").Value = 467
$ws.Range("D16").Value = 423
$ws.Range("F16").Value = 45566
$ws.Range("H16").Value = 45658
$ws.Range("C17").Value = 383
$ws.Range("D17").Value = 407
$ws.Range("F17").Value = 45597
$ws.Range("H17").Value = 45658
$ws.Range("D18").Value = 273
$ws.Range("E18").Value = 30407
$ws.Range("F18").Value = 45597
$ws.Range("H18").Value = 45658
$ws.Range("D19").Value = 411
$ws.Range("E19").Value = 30407
$ws.Range("F19").Value = 45597
$ws.Range("H19").Value = 45658
$ws.Range("C20").Value = 493
$ws.Range("F20").Value = 45566
$ws.Range("G20").Value = 30468
$ws.Range("H20").Value = 45658
$ws.Range("C21").Value = 323
$ws.Range("F21").Value = 45597
$ws.Range("G21").Value = 30468
$ws.Range("H21").Value = 45658
$ws.Range("C22").Value = 336
$ws.Range("D22").Value = 380
$ws.Range("F22").Value = 45597
$ws.Range("H22").Value = 45658
$ws.Range("D23").Value = 328
$ws.Range("E23").Value = 30376
$ws.Range("F23").Value = 45566
$ws.Range("H23").Value = 45658
$ws.Range("C24").Value = 345
$ws.Range("D24").Value = 326
$ws.Range("F24").Value = 45566
$ws.Range("H24").Value = 45658
$ws.Range("G25").Value = 30468
$ws.Range("H25").Value = 45658
$ws.Range("D26").Value = 393
$ws.Range("E26").Value = 30376
$ws.Range("F26").Value = 45566
$ws.Range("H26").Value = 45658
$ws.Range("D27").Value = 230
$ws.Range("E27").Value = 30407
$ws.Range("F27").Value = 45597
$ws.Range("H27").Value = 45658
$ws.Range("C28").Value = 406
$ws.Range("F28").Value = 45566
$ws.Range("G28").Value = 30468
$ws.Range("H28").Value = 45658
$ws.Range("E29").Value = 30407
$ws.Range("F29").Value = 45597
$ws.Range("G29").Value = 30468
$ws.Range("H29").Value = 45658
$ws.Range("D30").Value = 411
$ws.Range("H30").Value = 45658
$ws.Range("C31").Value = 417
$ws.Range("D31").Value = 330
$ws.Range("F31").Value = 45566
$ws.Range("H31").Value = 45658
$ws.Range("D32").Value = 423
$ws.Range("E32").Value = 30407
$ws.Range("F32").Value = 45597
$ws.Range("H32").Value = 45658
$ws.Range("C33").Value = 479
$ws.Range("D33").Value = 330
$ws.Range("F33").Value = 45597
$ws.Range("H33").Value = 45658
$ws.Range("C34").Value = 372
$ws.Range("D34").Value = 380
$ws.Range("F34").Value = 45597
$ws.Range("H34").Value = 45658
